$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.620.78'
$ws.Range("E2").Value = '  +3.80%  '

$ws.Range("D3").Value = '2.434.72'
$ws.Range("E3").Value = '  +2.56%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '312.12'
$ws.Range("E5").Value = '  +3.40%  '

$ws.Range("D6").Value = '102.04'
$ws.Range("E6").Value = '  +6.08%  '

$ws.Range("D7").Value = '0.515'
$ws.Range("E7").Value = '  +2.02%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '0.510'
$ws.Range("E9").Value = '  +2.49%  '

$ws.Range("D10").Value = '35.47'
$ws.Range("E10").Value = '  +3.87%  '

$ws.Range("D11").Value = '0.0802'
$ws.Range("E11").Value = '  +1.74%  '

$ws.Range("E12").Value = '  +0.97%  '

$ws.Range("D13").Value = '18.85'
$ws.Range("E13").Value = '  +3.76%  '

$ws.Range("E14").Value = '  +3.05%  '

$ws.Range("D15").Value = '2.813.09'
$ws.Range("E15").Value = '  +2.40%  '

$ws.Range("D16").Value = '2.411.07'
$ws.Range("E16").Value = '  -0.62%  '

$ws.Range("D17").Value = '0.839'
$ws.Range("E17").Value = '  +4.62%  '

$ws.Range("D18").Value = '44.510.87'
$ws.Range("E18").Value = '  +3.60%  '

$ws.Range("E19").Value = '  +2.76%  '

$ws.Range("D20").Value = '6.43'
$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("D21").Value = '0.0₃0911'
$ws.Range("E21").Value = '  +2.60%  '

$ws.Range("D22").Value = '68.95'
$ws.Range("E22").Value = '  +1.18%  '

$ws.Range("E23").Value = '  +4.34%  '

$ws.Range("D24").Value = '241.44'
$ws.Range("E24").Value = '  +2.76%  '

$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +1.93%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").Value = '25.31'
$ws.Range("E27").Value = '  +1.76%  '

$ws.Range("D28").Value = '2.27'
$ws.Range("E28").Value = '  -3.67%  '

$ws.Range("E29").Value = '  +4.96%  '

$ws.Range("D30").Value = '33.50'
$ws.Range("E30").Value = '  +6.32%  '

$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '0.123'
$ws.Range("E31").Value = '  +16.80%  '

$ws.Range("B32").Value = 'Celestia'
$ws.Range("C32").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D32").Value = '19.50'
$ws.Range("E32").Value = '  +11.29%  '

$ws.Range("D33").Value = '5.21'
$ws.Range("E33").Value = '  +2.76%  '

$ws.Range("E34").Value = '  +0.19%  '

$ws.Range("D35").Value = '0.0762'
$ws.Range("E35").Value = '  +4.06%  '

$ws.Range("D36").Value = '1.92'
$ws.Range("E36").Value = '  +3.02%  '

$ws.Range("D37").Value = '4.52'
$ws.Range("E37").Value = '  +3.57%  '

$ws.Range("D38").Value = '2.92'
$ws.Range("E38").Value = '  +4.28%  '

$ws.Range("D39").Value = '126.65'
$ws.Range("E39").Value = '  +8.30%  '

$ws.Range("D40").Value = '2.24'
$ws.Range("E40").Value = '  -3.01%  '

$ws.Range("E41").Value = '  +1.15%  '

$ws.Range("D42").Value = '21.75'
$ws.Range("E42").Value = '  -1.15%  '

$ws.Range("D43").Value = '0.0291'
$ws.Range("E43").Value = '  +3.62%  '

$ws.Range("D44").Value = '1.950.34'
$ws.Range("E44").Value = '  +0.27%  '

$ws.Range("E45").Value = '  +2.21%  '

$ws.Range("D46").Value = '2.97'
$ws.Range("E46").Value = '  +8.96%  '

$ws.Range("D47").Value = '9.73'
$ws.Range("E47").Value = '  +6.34%  '

$ws.Range("E48").Value = '  +11.18%  '

$ws.Range("D49").Value = '2.676.04'
$ws.Range("E49").Value = '  +2.75%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '53.47'
$ws.Range("E50").Value = '  +2.86%  '

$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = '74.00'
$ws.Range("E51").Value = '  +2.58%  '
